# Edit script implementing the "Messreihen" diff:
#  - sharedStrings: add new M08-M13 / D06-D13 rows of data + descriptions,
#    remove obsolete D01/D02 placeholders (superseded by renumbered D06/D07)
#  - sheet1: add "Versorgung" column (F), extend data table through row 39,
#    adjust view (selection) and column F width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 17
Set-Cell "B17" "R1 / Ω"
Set-Cell "C17" "R2 / Ω"
Set-Cell "D17" "Fs / Hz"
Set-Cell "E17" "Verstärker"
Set-Cell "F17" "Versorgung"
Set-Cell "G17" "Beschreibung"

# Row 18
Set-Cell "A18" "M01"
Set-Cell "B18" "1M"
Set-Cell "C18" "33K"
Set-Cell "D18" 1000
Set-Cell "E18" "ohne"
Set-Cell "F18" "Schaltnetzteil"
Set-Cell "G18" "I um 180° gedreht da falsch angesteckt, Spannungsteiler 1M 33K"

# Row 19
Set-Cell "A19" "M02"
Set-Cell "B19" "1M"
Set-Cell "C19" "33K"
Set-Cell "D19" 1000
Set-Cell "E19" "ohne"
Set-Cell "F19" "Schaltnetzteil"
Set-Cell "G19" "Messung mit richtiger Polung wiederholt, Spannungsteiler 1M 33K"

# Row 20
Set-Cell "A20" "M03"
Set-Cell "B20" "1M"
Set-Cell "C20" "33K"
Set-Cell "D20" 1000
Set-Cell "E20" "ohne"
Set-Cell "F20" "Schaltnetzteil"
Set-Cell "G20" "Beide Eingänge messen U, Spannungsteiler 1M 33K, Phasenverschiebung kann nur vom zeitlichen Versatz der Abtastung kommen"

# Row 21
Set-Cell "A21" "M04"
Set-Cell "B21" "18K"
Set-Cell "C21" 570
Set-Cell "D21" 1000
Set-Cell "E21" "ohne"
Set-Cell "F21" "Schaltnetzteil"
Set-Cell "G21" "Neue Messreihe R1 alternierend, Spannungsteiler 18K 570"

# Row 22
Set-Cell "A22" "M05"
Set-Cell "B22" "18K"
Set-Cell "C22" 570
Set-Cell "D22" 500
Set-Cell "E22" "ohne"
Set-Cell "F22" "Schaltnetzteil"
Set-Cell "G22" "Wie M04 nur andere Fs"

# Row 23
Set-Cell "A23" "M06"
Set-Cell "B23" "18K"
Set-Cell "C23" 570
Set-Cell "D23" 500
Set-Cell "E23" "ohne"
Set-Cell "F23" "Schaltnetzteil"
Set-Cell "G23" "I wieder zurück gedreht, Shunt war überbrückt, jetzt mehr Spannung an ADC1, aber 180° verdreht"

# Row 24
Set-Cell "A24" "M07"
Set-Cell "B24" "18K"
Set-Cell "C24" 570
Set-Cell "D24" 500
Set-Cell "E24" "ohne"
Set-Cell "F24" "Schaltnetzteil"
Set-Cell "G24" "i_temp wird nun im Porgamm negativ berechnet "

# Row 25
Set-Cell "A25" "M08"
Set-Cell "B25" "18K"
Set-Cell "C25" 570
Set-Cell "D25" 500
Set-Cell "E25" "ohne"
Set-Cell "F25" "Schaltnetzteil"
Set-Cell "G25" "gleich M07"

# Row 26
Set-Cell "A26" "M09"
Set-Cell "B26" "18K"
Set-Cell "C26" 570
Set-Cell "D26" 500
Set-Cell "E26" 4.7
Set-Cell "F26" "Schaltnetzteil"
Set-Cell "G26" "nicht-invertierender Verstärker und -i_temp, Übersteuern von i ab Messung 11!"

# Row 27
Set-Cell "A27" "M10"
Set-Cell "B27" "18K"
Set-Cell "C27" 570
Set-Cell "D27" 500
Set-Cell "E27" -4.7
Set-Cell "F27" "Schaltnetzteil"
Set-Cell "G27" "invertierender Verstärker, Übersteuern von i bei Messung 12"

# Row 28
Set-Cell "A28" "M11"
Set-Cell "B28" "18K"
Set-Cell "C28" 570
Set-Cell "D28" 500
Set-Cell "E28" -4.7
Set-Cell "F28" 7805
Set-Cell "G28" "invertierender Verstärker"

# Row 29
Set-Cell "A29" "M12"
Set-Cell "B29" "18K"
Set-Cell "C29" 570
Set-Cell "D29" 500
Set-Cell "E29" 4.7
Set-Cell "F29" 7805
Set-Cell "G29" "nicht-invertierender Verstärker und -i_temp, Übersteuern von i ab Messung 11!"

# Row 30
Set-Cell "A30" "M13"
Set-Cell "B30" "18K"
Set-Cell "C30" 570
Set-Cell "D30" 500
Set-Cell "E30" "ohne"
Set-Cell "F30" 7805
Set-Cell "G30" "ohne Versträker"

# Row 32
Set-Cell "A32" "D06"
Set-Cell "B32" "18K"
Set-Cell "C32" 570
Set-Cell "D32" 500
Set-Cell "E32" "ohne"
Set-Cell "F32" "Schaltnetzteil"
Set-Cell "G32" "Messung und Berechnung durch den µC, Übertragung der Leistungsdaten"

# Row 33
Set-Cell "A33" "D07"
Set-Cell "B33" "18K"
Set-Cell "C33" 570
Set-Cell "D33" 500
Set-Cell "E33" "ohne"
Set-Cell "F33" "Schaltnetzteil"
Set-Cell "G33" "Berechnete Leisuntgsdaten mit vorher gedrehtem Vorzeichen von i_temp"

# Row 34
Set-Cell "A34" "D08"
Set-Cell "B34" "18K"
Set-Cell "C34" 570
Set-Cell "D34" 500
Set-Cell "E34" "ohne"
Set-Cell "F34" "Schaltnetzteil"
Set-Cell "G34" "Kein Verstärker mit drehung VZ i_temp für positive P "

# Row 35
Set-Cell "A35" "D09"
Set-Cell "B35" "18K"
Set-Cell "C35" 570
Set-Cell "D35" 500
Set-Cell "E35" 4.7
Set-Cell "F35" "Schaltnetzteil"
Set-Cell "G35" "nicht-invertierender Verstärker und -i_temp, Übersteuern von i!"

# Row 36
Set-Cell "A36" "D10"
Set-Cell "B36" "18K"
Set-Cell "C36" 570
Set-Cell "D36" 500
Set-Cell "E36" -4.7
Set-Cell "F36" "Schaltnetzteil"
Set-Cell "G36" "invertierender Verstärker, Übersteuern von i bei Messung 12"

# Row 37
Set-Cell "A37" "D11"
Set-Cell "B37" "18K"
Set-Cell "C37" 570
Set-Cell "D37" 500
Set-Cell "E37" -4.7
Set-Cell "F37" 7805
Set-Cell "G37" "invertierender Verstärker"

# Row 38
Set-Cell "A38" "D12"
Set-Cell "B38" "18K"
Set-Cell "C38" 570
Set-Cell "D38" 500
Set-Cell "E38" 4.7
Set-Cell "F38" 7805
Set-Cell "G38" "nicht-invertierender Verstärker und -i_temp, Übersteuern von i! "

# Row 39
Set-Cell "A39" "D13"
Set-Cell "B39" "18K"
Set-Cell "C39" 570
Set-Cell "D39" 500
Set-Cell "E39" "ohne"
Set-Cell "F39" 7805
Set-Cell "G39" "ohne Versträker"

# Widen column F ("Versorgung") to fit its new content
$ws.Columns.Item(6).ColumnWidth = 12.5

# Update the view: scroll so row 4 is at the top and select G39 (last edited cell)
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G39").Select()

